$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 82.40641713471288
$ws.Range("B3").Value = 86.73939412659749
$ws.Range("B4").Value = 89.80045417174072
$ws.Range("H5").Value = 95.88925451112898
$ws.Range("H6").Value = 95.88701720109172
$ws.Range("H7").Value = 95.87308791677911
$ws.Range("C8").Value = 98.41432840090805
$ws.Range("C9").Value = 97.0255619905421
$ws.Range("C10").Value = 98.11787034460517
$ws.Range("D11").Value = 99.32639786170625
$ws.Range("D12").Value = 99.21500274414939
$ws.Range("D13").Value = 99.27432671514555
$ws.Range("E14").Value = 98.72743293293674
$ws.Range("E15").Value = 98.79197456695178
$ws.Range("E16").Value = 98.7566203622452
$ws.Range("F17").Value = 98.2289834120438
$ws.Range("F18").Value = 98.29785025341207
$ws.Range("F19").Value = 98.21321661171832
$ws.Range("G20").Value = 97.25361541597267
$ws.Range("G21").Value = 97.33924396059854
$ws.Range("G22").Value = 97.28763606178265
$ws.Range("B23").Value = 88.85368895367749
$ws.Range("B24").Value = 91.8468054008783
$ws.Range("H25").Value = 95.91481121943517
$ws.Range("H26").Value = 95.87691392297177
$ws.Range("C27").Value = 97.7092903514183
$ws.Range("C28").Value = 97.88577715113331
$ws.Range("D29").Value = 99.28963624262335
$ws.Range("D30").Value = 99.26375222298115
$ws.Range("E31").Value = 98.78417761786724
$ws.Range("E32").Value = 98.74528000914586
$ws.Range("F33").Value = 98.25490596901302
$ws.Range("F34").Value = 98.32839331878446
$ws.Range("G35").Value = 97.22558774092572
$ws.Range("G36").Value = 97.24699659693512
$ws.Range("B37").Value = 90.53224334979524
$ws.Range("B38").Value = 92.61237230484846
$ws.Range("H39").Value = 95.8771402532327
$ws.Range("H40").Value = 95.81760157069749
$ws.Range("C41").Value = 98.46355862844388
$ws.Range("C42").Value = 98.11542575674849
$ws.Range("D43").Value = 99.2903186757181
$ws.Range("D44").Value = 99.23788476988348
$ws.Range("E45").Value = 98.73713926589632
$ws.Range("E46").Value = 98.67393416462255
$ws.Range("F47").Value = 98.29375315960617
$ws.Range("F48").Value = 98.16132394066372
$ws.Range("G49").Value = 97.23947772932038
$ws.Range("G50").Value = 97.34820116697385
